$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.303.34'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.017.45'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.60'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.44'
$ws.Range("E6").Value = '  +2.63%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.017.54'
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("E10").Value = '  +11.07%  '
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.49'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("E15").Value = '  +2.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.518.16'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.255.82'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.023.19'
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.31'
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.22'
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.43'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.32'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("E25").Value = '  +2.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.82'
$ws.Range("E26").Value = '  +11.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.01'
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  +2.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.18'
$ws.Range("E31").Value = '  +3.91%  '
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("E34").Value = '  +1.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0850'
$ws.Range("E35").Value = '  +6.23%  '
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("E39").Value = '  +5.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.12'
$ws.Range("E40").Value = '  -0.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.04'
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.50'
$ws.Range("E43").Value = '  +11.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.284'
$ws.Range("E44").Value = '  +6.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '392.19'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.742.43'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.77'
$ws.Range("E48").Value = '  +3.98%  '
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("E51").Value = '  -1.08%  '
